$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2663.3333
$ws.Range("I62").Value = 1992
$ws.Range("J62").Value = 2999
$ws.Range("K62").Value = 1992
$ws.Range("L62").Value = 2999
$ws.Range("M62").Value = -1368
$ws.Range("H65").Value = 2663.3333
$ws.Range("I65").Value = 1992
$ws.Range("J65").Value = 2999
$ws.Range("K65").Value = 9960
$ws.Range("L65").Value = 14995
$ws.Range("M65").Value = -6840
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1000
$ws.Range("N107").Value = -4840
$ws.Range("M107").ClearContents()
$ws.Range("H111").Value = 25001724
$ws.Range("I111").Value = 50000250
$ws.Range("J111").Value = 3200
$ws.Range("K111").Value = 150000750
$ws.Range("L111").Value = 9600
$ws.Range("M111").Value = -149997683
$ws.Range("N111").Value = -15734
$ws.Range("H132").Value = 1073.5143
$ws.Range("I132").Value = 987.0606
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 2961.1818
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -431.1818000000003
$ws.Range("H138").Value = 3822.6
$ws.Range("I138").Value = 3873.125
$ws.Range("J138").Value = 3764.8572
$ws.Range("K138").Value = 11619.375
$ws.Range("L138").Value = 11294.5716
$ws.Range("M138").Value = -6479.375
$ws.Range("N138").Value = -21574.5716
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4590.95
$ws.Range("I61").Value = 2880
$ws.Range("J61").Value = 8583.166999999999
$ws.Range("K61").Value = 2880
$ws.Range("L61").Value = 8583.166999999999
$ws.Range("M61").Value = -2668
$ws.Range("H63").Value = 5383
$ws.Range("I63").Value = 5418.1665
$ws.Range("J63").Value = 5277.5
$ws.Range("K63").Value = 5418.1665
$ws.Range("L63").Value = 5277.5
$ws.Range("M63").Value = -4732.1665
$ws.Range("N63").Value = -6649.5
$ws.Range("H66").Value = 5383
$ws.Range("I66").Value = 5418.1665
$ws.Range("J66").Value = 5277.5
$ws.Range("K66").Value = 27090.8325
$ws.Range("L66").Value = 26387.5
$ws.Range("M66").Value = -23658.8325
$ws.Range("N66").Value = -33251.5
$ws.Range("H74").Value = 1297.6666
$ws.Range("I74").Value = 847.4474
$ws.Range("J74").Value = 5574.75
$ws.Range("K74").Value = 847.4474
$ws.Range("L74").Value = 5574.75
$ws.Range("M74").Value = 26.55259999999998
$ws.Range("N74").Value = -7322.75
$ws.Range("H77").Value = 1297.6666
$ws.Range("I77").Value = 847.4474
$ws.Range("J77").Value = 5574.75
$ws.Range("K77").Value = 4237.237
$ws.Range("L77").Value = 27873.75
$ws.Range("M77").Value = 130.7629999999999
$ws.Range("N77").Value = -36609.75
$ws.Range("H110").Value = 2993
$ws.Range("I110").Value = 1782.2
$ws.Range("J110").Value = 4506.5
$ws.Range("K110").Value = 1782.2
$ws.Range("L110").Value = 4506.5
$ws.Range("M110").Value = 262.8
$ws.Range("H132").Value = 1892.5714
$ws.Range("I132").Value = 1130.2727
$ws.Range("J132").Value = 2731.1
$ws.Range("K132").Value = 3390.8181
$ws.Range("L132").Value = 8193.299999999999
$ws.Range("M132").Value = -860.8181
$ws.Range("H136").Value = 4590.95
$ws.Range("I136").Value = 2880
$ws.Range("J136").Value = 8583.166999999999
$ws.Range("K136").Value = 8640
$ws.Range("L136").Value = 25749.501
$ws.Range("M136").Value = -6090
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 326.1
$ws.Range("I94").Value = 332.17242
$ws.Range("J94").Value = 150
$ws.Range("K94").Value = 332.17242
$ws.Range("L94").Value = 150
$ws.Range("M94").Value = 118.82758
$ws.Range("I99").Value = 1522
$ws.Range("J99").Value = 2333.3333
$ws.Range("K99").Value = 1522
$ws.Range("L99").Value = 2333.3333
$ws.Range("M99").Value = -24
$ws.Range("N99").Value = -5329.3333
$ws.Range("H107").Value = 3044.2856
$ws.Range("I107").Value = 3044.2856
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3044.2856
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1124.2856
$ws.Range("H134").Value = 15506.333
$ws.Range("I134").Value = 15695.875
$ws.Range("J134").Value = 14899.8
$ws.Range("K134").Value = 47087.625
$ws.Range("L134").Value = 44699.39999999999
$ws.Range("M134").Value = -44552.625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 679.6
$ws.Range("I22").Value = 349.5
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 349.5
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 0.5
$ws.Range("H62").Value = 2501.9
$ws.Range("I62").Value = 2311.6667
$ws.Range("J62").Value = 2787.25
$ws.Range("K62").Value = 2311.6667
$ws.Range("L62").Value = 2787.25
$ws.Range("M62").Value = -1687.6667
$ws.Range("N62").Value = -4035.25
$ws.Range("H65").Value = 2501.9
$ws.Range("I65").Value = 2311.6667
$ws.Range("J65").Value = 2787.25
$ws.Range("K65").Value = 11558.3335
$ws.Range("L65").Value = 13936.25
$ws.Range("M65").Value = -8438.333500000001
$ws.Range("N65").Value = -20176.25
$ws.Range("H107").Value = 299.47827
$ws.Range("I107").Value = 244.4
$ws.Range("J107").Value = 666.6667
$ws.Range("K107").Value = 244.4
$ws.Range("L107").Value = 666.6667
$ws.Range("M107").Value = 1675.6
$ws.Range("H132").Value = 2026.871
$ws.Range("I132").Value = 1139.8235
$ws.Range("J132").Value = 3104
$ws.Range("K132").Value = 3419.4705
$ws.Range("L132").Value = 9312
$ws.Range("M132").Value = -889.4704999999999
$ws.Range("H134").Value = 1049.6111
$ws.Range("I134").Value = 929.0625
$ws.Range("J134").Value = 2014
$ws.Range("K134").Value = 2787.1875
$ws.Range("L134").Value = 6042
$ws.Range("M134").Value = -252.1875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 181
$ws.Range("I10").Value = 181
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 543
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -404
$ws.Range("H92").Value = 300
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 300
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 900
$ws.Range("M92").ClearContents()
$ws.Range("H98").Value = 933.3333
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 933.3333
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 2799.9999
$ws.Range("N98").Value = -5795.9999
$ws.Range("M98").ClearContents()
$ws.Range("H107").Value = 843.1177
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 877.0625
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 2631.1875
$ws.Range("M107").Value = 1020
$ws.Range("N107").Value = -6471.1875
$ws.Range("H112").Value = 999
$ws.Range("I112").Value = 999
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 2997
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -1889
$ws.Range("N112").ClearContents()
$ws.Range("H125").Value = 5953.3335
$ws.Range("I125").Value = 1430
$ws.Range("J125").Value = 15000
$ws.Range("K125").Value = 4290
$ws.Range("L125").Value = 45000
$ws.Range("M125").Value = 630
$ws.Range("N125").Value = -54840
$ws.Range("H131").Value = 746.49
$ws.Range("I131").Value = 462.4
$ws.Range("J131").Value = 778.05554
$ws.Range("K131").Value = 1387.2
$ws.Range("L131").Value = 2334.16662
$ws.Range("M131").Value = 3652.8
$ws.Range("N131").Value = -12414.16662
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4796
$ws.Range("I102").Value = 5244.75
$ws.Range("J102").Value = 3001
$ws.Range("K102").Value = 5244.75
$ws.Range("L102").Value = 3001
$ws.Range("M102").Value = -3622.75
$ws.Range("H113").Value = 1224.5
$ws.Range("I113").Value = 999
$ws.Range("J113").Value = 1299.6666
$ws.Range("K113").Value = 999
$ws.Range("L113").Value = 1299.6666
$ws.Range("M113").Value = 1171
$ws.Range("N113").Value = -5639.6666
$ws.Range("H132").Value = 3370
$ws.Range("I132").Value = 2872.1052
$ws.Range("J132").Value = 8100
$ws.Range("K132").Value = 8616.3156
$ws.Range("L132").Value = 24300
$ws.Range("M132").Value = -6086.3156
$ws.Range("N132").Value = -29360
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2214.8948
$ws.Range("I132").Value = 2091.25
$ws.Range("J132").Value = 2304.818
$ws.Range("K132").Value = 6273.75
$ws.Range("L132").Value = 6914.454000000001
$ws.Range("M132").Value = -3743.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 17500
$ws.Range("I39").Value = 5000
$ws.Range("J39").Value = 30000
$ws.Range("K39").Value = 5000
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = -4587
$ws.Range("N39").Value = -30826
$ws.Range("H100").Value = 684.6667
$ws.Range("I100").Value = 449
$ws.Range("J100").Value = 873.2
$ws.Range("K100").Value = 898
$ws.Range("L100").Value = 1746.4
$ws.Range("M100").Value = -357
$ws.Range("N100").Value = -2828.4
$ws.Range("H122").Value = 79048.10000000001
$ws.Range("I122").Value = 112025.14
$ws.Range("J122").Value = 2101.6667
$ws.Range("K122").Value = 336075.42
$ws.Range("L122").Value = 6305.000100000001
$ws.Range("M122").Value = -333625.42
$ws.Range("N122").Value = -11205.0001
